$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "60.752.50"
Set-TextValue "E2" "  +4.01%  "
Set-TextValue "D3" "2.654.82"
Set-TextValue "E3" "  +1.36%  "
Set-TextValue "E4" "  +0.06%  "
Set-TextValue "D5" "569.73"
Set-TextValue "E5" "  +6.92%  "
Set-TextValue "D6" "147.66"
Set-TextValue "E6" "  +3.91%  "
Set-TextValue "E7" "  -0.43%  "
Set-TextValue "E8" "  +6.99%  "
Set-TextValue "D9" "6.87"
Set-TextValue "E9" "  -0.68%  "
Set-TextValue "E10" "  +5.03%  "
Set-TextValue "D11" "0.143"
Set-TextValue "E11" "  +6.58%  "
Set-TextValue "E12" "  +3.57%  "
Set-TextValue "D13" "3.123.85"
Set-TextValue "E13" "  +1.17%  "
Set-TextValue "D14" "60.705.01"
Set-TextValue "E14" "  +4.06%  "
Set-TextValue "D15" "21.86"
Set-TextValue "E15" "  +5.93%  "
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "2.680.98"
Set-TextValue "E16" "  +2.28%  "
Set-TextValue "B17" "ShibaInu"
Set-TextValue "C17" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.0000138"
Set-TextValue "E17" "  +5.06%  "
Set-TextValue "E18" "  +4.31%  "
Set-TextValue "D19" "345.90"
Set-TextValue "E19" "  +3.66%  "
Set-TextValue "D20" "10.48"
Set-TextValue "E20" "  +3.87%  "
Set-TextValue "D21" "6.44"
Set-TextValue "E21" "  +3.62%  "
Set-TextValue "E22" "  +1.36%  "
Set-TextValue "E23" "  -0.21%  "
Set-TextValue "D24" "66.78"
Set-TextValue "E24" "  +0.68%  "
Set-TextValue "D25" "0.444"
Set-TextValue "E25" "  +7.27%  "
Set-TextValue "D26" "0.166"
Set-TextValue "E26" "  +2.35%  "
Set-TextValue "D27" "0.996"
Set-TextValue "E27" "  -0.38%  "
Set-TextValue "D28" "7.40"
Set-TextValue "E28" "  +4.76%  "
Set-TextValue "D29" "0.0₃0792"
Set-TextValue "E29" "  +8.44%  "
Set-TextValue "E30" "  -0.14%  "
Set-TextValue "E31" "  +5.12%  "
Set-TextValue "D32" "6.17"
Set-TextValue "E32" "  +5.66%  "
Set-TextValue "D33" "19.34"
Set-TextValue "D34" "155.24"
Set-TextValue "D35" "4.12"
Set-TextValue "E35" "  +6.61%  "
Set-TextValue "D36" "0.919"
Set-TextValue "E36" "  +8.57%  "
Set-TextValue "E37" "  +8.69%  "
Set-TextValue "E38" "  +14.32%  "
Set-TextValue "D39" "37.66"
Set-TextValue "E39" "  +1.35%  "
Set-TextValue "D40" "1.52"
Set-TextValue "E40" "  +8.35%  "
Set-TextValue "D41" "309.80"
Set-TextValue "E41" "  +11.00%  "
Set-TextValue "D42" "3.68"
Set-TextValue "E42" "  +3.62%  "
Set-TextValue "D43" "0.611"
Set-TextValue "E43" "  +3.02%  "
Set-TextValue "E44" "  -0.58%  "
Set-TextValue "D45" "0.0983"
Set-TextValue "E45" "  +5.27%  "
Set-TextValue "D46" "0.0552"
Set-TextValue "E46" "  +4.65%  "
Set-TextValue "D47" "19.64"
Set-TextValue "E47" "  +3.99%  "
Set-TextValue "E48" "  -0.04%  "
Set-TextValue "B49" "Aave"
Set-TextValue "C49" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D49" "126.02"
Set-TextValue "E49" "  +11.97%  "
Set-TextValue "B50" "VeChain"
Set-TextValue "C50" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D50" "0.0235"
Set-TextValue "E50" "  +5.61%  "
Set-TextValue "D51" "4.77"
Set-TextValue "E51" "  +8.00%  "
